# Auto-generated edit script: updates cryptos price/volume columns
# per the commit diff (GitHub Actions cryptos-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Writes $text as a literal Text-typed cell value. Excel (and the
    # COM-interop engine) auto-converts plain numeric-looking strings
    # (e.g. "206.18") into Number cells when assigned via .Value, which
    # would silently change the cell type from the Text it is in the
    # source workbook. Prefixing with a leading apostrophe forces Excel
    # to store the literal text, matching the original inline-string cells.
    $range.Value = "'" + $text
}

$ws.Range("D2").Value = "26.890.05"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.551.33"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.39%  "
Set-TextValue $ws.Range("D5") "206.18"
$ws.Range("E5").Value = "  +0.80%  "
Set-TextValue $ws.Range("D6") "0.484"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.40%  "
Set-TextValue $ws.Range("D9") "21.48"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  +0.30%  "
Set-TextValue $ws.Range("D11") "0.0857"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.772.64"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.577.89"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("E14").Value = "  +0.95%  "
Set-TextValue $ws.Range("D15") "0.513"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "26.898.85"
$ws.Range("E16").Value = "  +0.42%  "
Set-TextValue $ws.Range("D17") "61.55"
$ws.Range("E17").Value = "  +0.98%  "
Set-TextValue $ws.Range("D18") "213.82"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  -0.84%  "
Set-TextValue $ws.Range("D23") "9.16"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  -1.74%  "
Set-TextValue $ws.Range("D25") "153.17"
$ws.Range("E25").Value = "  +0.62%  "
Set-TextValue $ws.Range("D27") "14.85"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  -0.21%  "
Set-TextValue $ws.Range("D32") "3.22"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "1.366.06"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  +3.29%  "
Set-TextValue $ws.Range("D36") "0.969"
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("E37").Value = "  +0.51%  "
Set-TextValue $ws.Range("D38") "0.0163"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +0.78%  "
Set-TextValue $ws.Range("D42") "0.985"
$ws.Range("E42").Value = "  -0.42%  "
Set-TextValue $ws.Range("D43") "5.51"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("E44").Value = "  +3.65%  "
Set-TextValue $ws.Range("D45") "63.52"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "1.685.40"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +0.31%  "
Set-TextValue $ws.Range("D49") "0.0507"
$ws.Range("E49").Value = "  -0.16%  "
Set-TextValue $ws.Range("D50") "0.0953"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  +0.43%  "
